$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.067.01"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "2.673.35"
$ws.Range("E3").Value = "  +2.49%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "572.13"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").Value = "145.09"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("E10").Value = "  +1.17%  "
$ws.Range("E11").Value = "  +5.34%  "
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "3.145.50"
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("E14").Value = "  +11.40%  "
$ws.Range("D15").Value = "61.064.02"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "2.667.52"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("E18").Value = "  +3.17%  "
$ws.Range("D19").Value = "4.78"
$ws.Range("E19").Value = "  +2.69%  "
$ws.Range("D20").Value = "351.60"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").Value = "0.533"
$ws.Range("E23").Value = "  +1.91%  "
$ws.Range("D24").Value = "64.01"
$ws.Range("E24").Value = "  +0.98%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.163"
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("E27").Value = "  +5.29%  "
$ws.Range("D28").Value = "1.98"
$ws.Range("E28").Value = "  +5.93%  "
$ws.Range("D29").Value = "0.0₃0817"
$ws.Range("E29").Value = "  +2.75%  "
$ws.Range("D30").Value = "6.90"
$ws.Range("E30").Value = "  +8.09%  "
$ws.Range("D31").Value = "0.998"
$ws.Range("D32").Value = "165.31"
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("E34").Value = "  +10.97%  "
$ws.Range("E35").Value = "  +6.07%  "
$ws.Range("E36").Value = "  +6.31%  "
$ws.Range("D37").Value = "1.67"
$ws.Range("E37").Value = "  +4.13%  "
$ws.Range("D38").Value = "335.98"
$ws.Range("E38").Value = "  +11.78%  "
$ws.Range("D39").Value = "4.04"
$ws.Range("E39").Value = "  +4.44%  "
$ws.Range("D40").Value = "38.60"
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("D41").Value = "0.886"
$ws.Range("E41").Value = "  +5.10%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "5.22"
$ws.Range("E42").Value = "  +4.80%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "20.48"
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("D44").Value = "134.65"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0564"
$ws.Range("E45").Value = "  +2.98%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "0.100"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("E47").Value = "  +3.24%  "
$ws.Range("D48").Value = "0.617"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").Value = "20.55"
$ws.Range("E49").Value = "  +3.52%  "
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").Value = "2.099.97"
$ws.Range("E51").Value = "  +3.62%  "
